$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(394).Insert()

$ws.Cells.Item(394, 1).Value = 5
$ws.Cells.Item(394, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(394, 3).Value = "Maule"
$ws.Cells.Item(394, 4).Value = 45229
$ws.Cells.Item(394, 5).Value = 7
$ws.Cells.Item(394, 6).Value = 100112008
$ws.Cells.Item(394, 7).Value = "Coliflor"
$ws.Cells.Item(394, 8).Value = "Sin especificar"
$ws.Cells.Item(394, 9).Value = "Primera"
$ws.Cells.Item(394, 10).Value = 2000
$ws.Cells.Item(394, 11).Value = 1000
$ws.Cells.Item(394, 12).Value = 1000
$ws.Cells.Item(394, 13).Value = 1000
$ws.Cells.Item(394, 14).Value = "$/unidad"
$ws.Cells.Item(394, 15).Value = "Región del Maule"
$ws.Cells.Item(394, 16).Value = 1000
$ws.Cells.Item(394, 17).Value = 1
$ws.Cells.Item(394, 18).Value = "Hortaliza"
